$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells we are about to write to be treated as literal text
# (values like "590.28" or "0.0000276" would otherwise be auto-coerced into numbers),
# then reset the style back to Normal so no residual formatting/style is left behind.
$dCells = @("D2","D3","D5","D6","D7","D8","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D39","D41","D42","D43","D44","D45","D46","D47","D48","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "69.011.73"
$ws.Range("D3").Value = "3.515.22"
$ws.Range("D5").Value = "590.28"
$ws.Range("D6").Value = "171.30"
$ws.Range("D7").Value = "0.614"
$ws.Range("D8").Value = "3.512.34"
$ws.Range("D10").Value = "0.189"
$ws.Range("D11").Value = "6.93"
$ws.Range("D12").Value = "0.581"
$ws.Range("D13").Value = "47.12"
$ws.Range("D14").Value = "0.0000276"
$ws.Range("D15").Value = "4.086.18"
$ws.Range("D16").Value = "8.47"
$ws.Range("D17").Value = "621.20"
$ws.Range("D18").Value = "3.524.81"
$ws.Range("D19").Value = "69.123.78"
$ws.Range("D20").Value = "0.121"
$ws.Range("D21").Value = "17.41"
$ws.Range("D22").Value = "11.14"
$ws.Range("D23").Value = "0.885"
$ws.Range("D24").Value = "15.88"
$ws.Range("D25").Value = "96.88"
$ws.Range("D28").Value = "2.63"
$ws.Range("D29").Value = "9.25"
$ws.Range("D30").Value = "32.68"
$ws.Range("D31").Value = "3.13"
$ws.Range("D32").Value = "8.50"
$ws.Range("D33").Value = "1.32"
$ws.Range("D34").Value = "6.92"
$ws.Range("D35").Value = "635.53"
$ws.Range("D36").Value = "10.77"
$ws.Range("D37").Value = "3.47"
$ws.Range("D39").Value = "57.33"
$ws.Range("D41").Value = "0.0454"
$ws.Range("D42").Value = "0.135"
$ws.Range("D43").Value = "3.377.05"
$ws.Range("D44").Value = "0.327"
$ws.Range("D45").Value = "32.84"
$ws.Range("D46").Value = "0.0₃0695"
$ws.Range("D47").Value = "2.53"
$ws.Range("D48").Value = "2.77"
$ws.Range("D50").Value = "133.04"
$ws.Range("D51").Value = "5.61"

foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}

# Remaining (non price-column) text updates
$ws.Range("E2").Value = "  -2.32%  "
$ws.Range("E3").Value = "  -3.40%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("E6").Value = "  -2.35%  "
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("E8").Value = "  -3.29%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  -4.11%  "
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("E12").Value = "  -4.35%  "
$ws.Range("E13").Value = "  -2.57%  "
$ws.Range("E14").Value = "  -2.64%  "
$ws.Range("E16").Value = "  -4.85%  "
$ws.Range("E17").Value = "  -7.90%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("E18").Value = "  -3.18%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("E19").Value = "  -2.22%  "
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("E21").Value = "  -2.20%  "
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("E23").Value = "  -5.94%  "
$ws.Range("E24").Value = "  -7.21%  "
$ws.Range("E25").Value = "  -3.09%  "
$ws.Range("E26").Value = "  -1.55%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -5.82%  "
$ws.Range("E29").Value = "  -6.58%  "
$ws.Range("E30").Value = "  -5.95%  "
$ws.Range("E31").Value = "  -5.63%  "
$ws.Range("E32").Value = "  -5.47%  "
$ws.Range("E33").Value = "  -4.79%  "
$ws.Range("E34").Value = "  -7.61%  "
$ws.Range("E35").Value = "  +8.80%  "
$ws.Range("E36").Value = "  -2.64%  "
$ws.Range("E37").Value = "  -12.55%  "
$ws.Range("E38").Value = "  -4.29%  "
$ws.Range("E39").Value = "  -1.67%  "
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("E43").Value = "  -5.25%  "
$ws.Range("E44").Value = "  -4.99%  "
$ws.Range("E45").Value = "  -4.99%  "
$ws.Range("E46").Value = "  -5.38%  "
$ws.Range("E47").Value = "  -5.55%  "
$ws.Range("E48").Value = "  -3.80%  "
$ws.Range("E49").Value = "  -2.67%  "
$ws.Range("E50").Value = "  -2.43%  "
$ws.Range("E51").Value = "  +12.97%  "
